$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 6708.013860684405
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6718.366062258094
